# Add 2022-Q4 data:
#  - insert a new "2022-Q4" sheet right after the "总计" (summary) sheet,
#    holding the fund-position detail rows for the new quarter
#  - insert a new row into "总计" summarizing that quarter and shift the
#    existing rows down (re-numbering the leading index column)

$wb = $excel.ActiveWorkbook

# --- 1. "总计" summary sheet: insert the 2022-Q4 row at the top of the data ---
$summary = $wb.Worksheets.Item(1)

$summary.Range("A2").EntireRow.Insert()

# EntireRow.Insert() drags the prior row's formatting down onto the new
# row; the other data rows carry no explicit style, so strip it back off.
$summary.Range("B2:D2").ClearFormats()

# carry over the index column's number/border style from the row below
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 1.09

# renumber the rest of the index column (rows shifted down by one)
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6

# --- 2. New "2022-Q4" detail sheet, placed right after "总计" ---
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

# Re-use the header / index-column formatting from an existing quarter sheet
# (same layout: bold bordered header row, bordered index column) so the new
# sheet matches the others instead of staying unstyled.
$template = $wb.Worksheets.Item(3)   # "2022-Q2" - same 8-column layout
$template.Range("A1:H3").Copy()
$q4.Range("A1").PasteSpecial(-4122)
$q4.Range("A3").Copy()
$q4.Range("A4").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# text-typed detail columns (codes / amounts kept as strings, matching the
# other quarter sheets) - format as text first so leading zeros survive
$q4.Range("B2:G4").NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "160106"
$q4.Range("C2").Value = "南方高增长混合（LOF）"
$q4.Range("D2").Value = "16.27"
$q4.Range("E2").Value = "88.50"
$q4.Range("F2").Value = "4.63"
$q4.Range("G2").Value = "0.7533"
$q4.Range("H2").Value = 9

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "160105"
$q4.Range("C3").Value = "南方积极配置混合（LOF）"
$q4.Range("D3").Value = "5.52"
$q4.Range("E3").Value = "89.92"
$q4.Range("F3").Value = "4.69"
$q4.Range("G3").Value = "0.2589"
$q4.Range("H3").Value = 10

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "000554"
$q4.Range("C4").Value = "南方中国梦灵活配置混合"
$q4.Range("D4").Value = "1.50"
$q4.Range("E4").Value = "92.96"
$q4.Range("F4").Value = "5.15"
$q4.Range("G4").Value = "0.0772"
$q4.Range("H4").Value = 8
